$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Price values in column D are plain-text cells in the source data
# (e.g. "42.00", "0.0910") even though they look numeric - a leading
# apostrophe forces Excel to keep them as text instead of silently
# coercing to a Number (which would strip meaningful trailing/leading
# zeros). Column E (percent change) is already safe as literal text
# because of the surrounding spaces.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'42.927.72"
$ws.Range("E2").Value = "  -0.09%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'2.208.54"
$ws.Range("E3").Value = "  -1.67%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.17%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'254.41"
$ws.Range("E5").Value = "  +3.86%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.615"
$ws.Range("E6").Value = "  -0.42%  "

# Row 7 - Solana
$ws.Range("D7").Value = "'76.07"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.01%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.592"
$ws.Range("E9").Value = "  -4.03%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'42.00"
$ws.Range("E10").Value = "  +2.88%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.0910"
$ws.Range("E11").Value = "  -2.46%  "

# Row 12 - now Polkadot (was TRON)
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'6.87"
$ws.Range("E12").Value = "  -1.17%  "

# Row 13 - now TRON (was Polkadot)
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.102"
$ws.Range("E13").Value = "  +0.73%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "'2.539.84"
$ws.Range("E14").Value = "  -0.81%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "'14.42"
$ws.Range("E15").Value = "  -1.36%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "'2.207.40"

# Row 17 - Polygon
$ws.Range("D17").Value = "'0.780"
$ws.Range("E17").Value = "  -3.25%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "'42.791.34"
$ws.Range("E18").Value = "  -0.17%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "'0.0000102"
$ws.Range("E19").Value = "  -2.22%  "

# Row 20 - Litecoin
$ws.Range("D20").Value = "'71.20"
$ws.Range("E20").Value = "  +0.02%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -1.06%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'229.30"
$ws.Range("E22").Value = "  -0.75%  "

# Row 23 - ImmutableX
$ws.Range("E23").Value = "  -1.05%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("E24").Value = "  -9.05%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "'10.63"
$ws.Range("E26").Value = "  -2.26%  "

# Row 27 - WEMIXToken
$ws.Range("E27").Value = "  -2.00%  "

# Row 28 - InjectiveProtocol
$ws.Range("D28").Value = "'39.26"
$ws.Range("E28").Value = "  +1.72%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "'2.26"
$ws.Range("E29").Value = "  +5.59%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -2.78%  "

# Row 31 - Monero
$ws.Range("D31").Value = "'173.15"
$ws.Range("E31").Value = "  -0.35%  "

# Row 32 - EthereumClassic
$ws.Range("E32").Value = "  -0.57%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.0856"
$ws.Range("E33").Value = "  +7.44%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "'5.19"
$ws.Range("E34").Value = "  -2.58%  "

# Row 35 - Stellar
$ws.Range("E35").Value = "  -1.10%  "

# Row 36 - Kaspa
$ws.Range("D36").Value = "'0.107"
$ws.Range("E36").Value = "  -1.87%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "'0.0354"
$ws.Range("E37").Value = "  +6.58%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  -0.55%  "

# Row 39 - Celestia
$ws.Range("E39").Value = "  -2.69%  "

# Row 40 - LidoDAOToken
$ws.Range("D40").Value = "'2.09"
$ws.Range("E40").Value = "  -1.81%  "

# Row 41 - NEARProtocol
$ws.Range("D41").Value = "'2.74"
$ws.Range("E41").Value = "  +17.47%  "

# Row 42 - Algorand
$ws.Range("E42").Value = "  -2.89%  "

# Row 43 - THORChain
$ws.Range("E43").Value = "  -5.26%  "

# Row 44 - MultiversX
$ws.Range("D44").Value = "'59.88"
$ws.Range("E44").Value = "  +0.09%  "

# Row 45 - Aave
$ws.Range("D45").Value = "'101.53"
$ws.Range("E45").Value = "  -4.72%  "

# Row 46 - Cronos
$ws.Range("D46").Value = "'0.0977"
$ws.Range("E46").Value = "  -1.59%  "

# Row 47 - FraxShare
$ws.Range("D47").Value = "'8.29"

# Row 48 - WOONetwork
$ws.Range("E48").Value = "  -1.75%  "

# Row 49 - ARBITRUM
$ws.Range("E49").Value = "  -0.33%  "

# Row 50 - TrustWalletToken
$ws.Range("E50").Value = "  -1.63%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "'2.433.31"
$ws.Range("E51").Value = "  -0.64%  "
